$wb = $excel.ActiveWorkbook

# --- "Entries" sheet: rename the audit entry and fix its notification count. ---
# (Updated first so the new shared-string "Audit denied other access and
# notify only" is registered before the Rules-sheet string that embeds it.)
$entries = $wb.Worksheets.Item("Entries")
$entries.Range("A6").Value = "Audit denied other access and notify only"
$entries.Range("J6").Value = 1

# --- "Rules" sheet: remove the "Grant full access to allowed full access USBs" row
# and update the remaining RO-access row's Entries reference. ---
$rules = $wb.Worksheets.Item("Rules")
$rules.Rows("3").Delete()
$rules.Range("E3").Value = "Deny other access, Audit denied other access and notify only"
$rules.Range("A3").Select()

# --- "Groups" sheet: move the saved cursor/selection. ---
$groups = $wb.Worksheets.Item("Groups")
$groups.Range("A4").Select()

# --- Make "Entries" the active sheet/selection, matching the saved view state. ---
$entries.Activate()
$entries.Range("A6").Select()
